$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.577.83"
$ws.Range("E2").Value = "  -2.62%  "

$ws.Range("D3").Value = "3.934.41"
$ws.Range("E3").Value = "  -2.57%  "

$ws.Range("E4").Value = "  -0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "535.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.19%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.22"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.23%  "

$ws.Range("D7").Value = "3.929.00"
$ws.Range("E7").Value = "  -2.57%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.685"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -6.67%  "

$ws.Range("E9").Value = "  +0.04%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.733"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.75%  "

$ws.Range("E11").Value = "  -5.46%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.13"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +15.16%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000316"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.36%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.53"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.48%  "

$ws.Range("D15").Value = "4.558.50"
$ws.Range("E15").Value = "  -2.78%  "

$ws.Range("D16").Value = "3.928.45"
$ws.Range("E16").Value = "  -3.20%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "20.26"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.04%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.66"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.15%  "

$ws.Range("E19").Value = "  -1.26%  "

$ws.Range("E20").Value = "  -4.10%  "

$ws.Range("D21").Value = "70.457.26"
$ws.Range("E21").Value = "  -2.63%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "421.77"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.78%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "95.91"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -8.57%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.53"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.75%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.15"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.25%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "14.11"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.34%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.33"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.05%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.52"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.89%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.87"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.21%  "

$ws.Range("E30").Value = "  +16.14%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "36.16"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.88%  "

$ws.Range("E32").Value = "  +11.19%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "49.23"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +14.92%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "13.20"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.25%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "679.55"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.63%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.128"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.90%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "63.89"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.42%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.433"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.99%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.46"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.27%  "

$ws.Range("B40").Value = "PEPE"
$ws.Range("C40").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D40").Value = "0.0₃0813"
$ws.Range("E40").Value = "  -4.95%  "

$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.149"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.39%  "

$ws.Range("E42").Value = "  +0.15%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.997"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.09%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.22"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.67%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0479"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.44%  "

$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.147"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -8.14%  "

$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.70"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.68%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.68"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.06%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.37"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.58%  "

$ws.Range("B50").Value = "FLOKI"
$ws.Range("C50").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000275"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.35%  "

$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.96"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.56%  "
